$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in K2 (Fecha:) from 02-03-2020 to 03-03-2020
# Force text storage so Excel doesn't auto-convert the string into a date serial number
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "03-03-2020"

# Clear the first student row (name + km expense) - Roldán Vara, Sergio / 5047 km
$ws.Range("A8").ClearContents()
$ws.Range("A8").WrapText = $false
$ws.Range("K8").ClearContents()
$ws.Range("K8").WrapText = $false

# Clear the second student row (name + km expense) - Torres Gijón, Beatriz / 6 km
$ws.Range("A9").ClearContents()
$ws.Range("A9").WrapText = $false
$ws.Range("K9").ClearContents()
$ws.Range("K9").WrapText = $false
